$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "48.033.50"
Set-TextValue "E2" "  -0.39%  "
Set-TextValue "D3" "2.486.11"
Set-TextValue "E3" "  -1.63%  "
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "317.32"
Set-TextValue "E5" "  -2.03%  "
Set-TextValue "D6" "105.49"
Set-TextValue "E6" "  -3.43%  "
Set-TextValue "D7" "0.519"
Set-TextValue "E7" "  -1.92%  "
Set-TextValue "E8" "  -0.07%  "
Set-TextValue "D9" "0.537"
Set-TextValue "E9" "  -3.60%  "
Set-TextValue "D10" "38.88"
Set-TextValue "E10" "  -5.29%  "
Set-TextValue "D11" "20.15"
Set-TextValue "E11" "  -1.57%  "
Set-TextValue "E12" "  -3.09%  "
Set-TextValue "E13" "  +0.11%  "
Set-TextValue "D14" "7.09"
Set-TextValue "E14" "  -2.72%  "
Set-TextValue "D15" "2.875.87"
Set-TextValue "E15" "  -1.59%  "
Set-TextValue "D16" "2.490.73"
Set-TextValue "E16" "  -1.59%  "
Set-TextValue "D17" "0.827"
Set-TextValue "E17" "  -3.82%  "
Set-TextValue "D18" "47.936.75"
Set-TextValue "E18" "  -0.24%  "
Set-TextValue "D19" "2.99"
Set-TextValue "E19" "  +10.34%  "
Set-TextValue "D20" "12.72"
Set-TextValue "E20" "  -4.43%  "
Set-TextValue "D21" "6.56"
Set-TextValue "E21" "  -1.30%  "
Set-TextValue "D22" "0.0₃0929"
Set-TextValue "E22" "  -2.22%  "
Set-TextValue "D23" "70.97"
Set-TextValue "E23" "  -1.70%  "
Set-TextValue "D24" "271.82"
Set-TextValue "E24" "  +0.85%  "
Set-TextValue "D25" "2.50"
Set-TextValue "E25" "  -3.34%  "
Set-TextValue "E26" "  +0.16%  "
Set-TextValue "D27" "25.67"
Set-TextValue "E27" "  -2.20%  "
Set-TextValue "D28" "2.28"
Set-TextValue "E28" "  +2.77%  "
Set-TextValue "D29" "9.70"
Set-TextValue "E29" "  -4.64%  "
Set-TextValue "E30" "  -4.13%  "
Set-TextValue "D31" "34.57"
Set-TextValue "E31" "  -3.64%  "
Set-TextValue "D32" "49.32"
Set-TextValue "E32" "  -0.60%  "
Set-TextValue "E33" "  -0.12%  "
Set-TextValue "D34" "19.01"
Set-TextValue "E34" "  -5.01%  "
Set-TextValue "D35" "5.26"
Set-TextValue "E35" "  -2.74%  "
Set-TextValue "D36" "0.0771"
Set-TextValue "E36" "  -3.06%  "
Set-TextValue "D37" "1.93"
Set-TextValue "E37" "  -3.15%  "
Set-TextValue "D38" "4.56"
Set-TextValue "E38" "  -3.77%  "
Set-TextValue "D39" "2.86"
Set-TextValue "E39" "  -5.02%  "
Set-TextValue "D40" "122.58"
Set-TextValue "E40" "  +2.08%  "
Set-TextValue "D41" "0.110"
Set-TextValue "E41" "  -2.02%  "
Set-TextValue "D43" "21.85"
Set-TextValue "E43" "  -1.41%  "
Set-TextValue "E44" "  +0.64%  "
Set-TextValue "D45" "2.001.57"
Set-TextValue "E45" "  -1.02%  "
Set-TextValue "D46" "3.17"
Set-TextValue "E46" "  -0.29%  "
Set-TextValue "E47" "  -0.84%  "
Set-TextValue "D48" "1.99"
Set-TextValue "E48" "  -2.32%  "
Set-TextValue "D49" "8.89"
Set-TextValue "E49" "  -2.97%  "
Set-TextValue "D50" "5.17"
Set-TextValue "E50" "  -1.95%  "
Set-TextValue "D51" "78.63"
Set-TextValue "E51" "  -1.48%  "
